{"js": "// Highlight quantitative impact metrics (percentages, dollar amounts, large\n// numbers) in bold + dark slate color (#2C3E50) across the resume body.\n//\n// For each target paragraph we locate it by a unique, stable marker\n// substring, then use the paragraph-scoped `search()` API to find each\n// metric token in left-to-right order and apply bold + color directly to\n// the returned sub-range. Word/Office.js automatically splits the\n// paragraph's run(s) around the highlighted sub-range.\n\nconst HIGHLIGHT_COLOR = \"#2C3E50\";\n\nasync function highlightTokens(context, paragraphMarker, tokens) {\n  const paragraphs = context.document.body.paragraphs;\n  paragraphs.load(\"items/text\");\n  await context.sync();\n\n  const target = paragraphs.items.find((p) => p.text.includes(paragraphMarker));\n  if (!target) {\n    throw new Error(`Paragraph not found for marker: ${paragraphMarker}`);\n  }\n\n  for (const token of tokens) {\n    const results = target.search(token, { matchCase: true });\n    results.load(\"items\");\n    await context.sync();\n\n    if (results.items.length === 0) {\n      throw new Error(`Token not found: ${token} (in paragraph: ${paragraphMarker})`);\n    }\n\n    const hit = results.items[0];\n    hit.font.bold = true;\n    hit.font.color = HIGHLIGHT_COLOR;\n    await context.sync();\n  }\n}\n\n// 1. Data Science & Political Analytics bullets (Siege Analytics)\nawait highlightTokens(context, \"Discovered systematic race coding errors\", [\"23%\", \"64%\"]);\nawait highlightTokens(context, \"Utilized advanced sampling methods\", [\n  \"\\u00B14.2%\",\n  \"\\u00B12.1%\",\n  \"71%\",\n  \"87%\",\n]);\nawait highlightTokens(context, \"Trigonometric algorithm for boundary estimation reduced\", [\n  \"73.5%\",\n  \"$4.7M\",\n]);\nawait highlightTokens(context, \"Built real-time FEC analysis systems\", [\"$2\"]);\n\n// 2. Data Products Manager (Helm/Murmuration) bullet\nawait highlightTokens(context, \"Modernized legacy ETL processes\", [\"57%\"]);\n\n// 3. KEY ACHIEVEMENTS AND IMPACT bullets\nawait highlightTokens(\n  context,\n  \"Algorithmic innovation: Pioneered trigonometric boundary estimation\",\n  [\"73.5%\"]\n);\nawait highlightTokens(context, \"savings enabled nonprofit access\", [\"$4.7M\"]);\nawait highlightTokens(context, \"Platform impact: Built redistricting system serving\", [\n  \"12,847\",\n]);\n", "ps1": "# Highlight quantitative impact metrics (percentages, dollar amounts, large\n# numbers) in bold + dark slate color (#2C3E50) across the resume body.\n#\n# For each target paragraph we locate it by a unique, stable substring, then\n# run a Find against the sub-range that starts where the previous highlight\n# ended (so repeated tokens such as \"73.5%\" or \"$4.7M\" are matched once, in\n# paragraph order, without jumping to a different paragraph that happens to\n# contain the same text).\n\nfunction Get-WdColor([string]$hex) {\n    $r = [Convert]::ToInt32($hex.Substring(0, 2), 16)\n    $g = [Convert]::ToInt32($hex.Substring(2, 2), 16)\n    $b = [Convert]::ToInt32($hex.Substring(4, 2), 16)\n    # Word/OLE colors are stored BGR, i.e. 0xBBGGRR.\n    return ($b * 65536) + ($g * 256) + $r\n}\n\n$HighlightColor = Get-WdColor('2C3E50')\n\nfunction Highlight-Tokens([string]$paragraphMarker, [string[]]$tokens) {\n    $d = $word.ActiveDocument\n    $target = $null\n    foreach ($p in $d.Paragraphs) {\n        if ($p.Range.Text.Contains($paragraphMarker)) {\n            $target = $p\n            break\n        }\n    }\n    if ($null -eq $target) {\n        throw \"Paragraph not found for marker: $paragraphMarker\"\n    }\n\n    $pEnd = $target.Range.End\n    $searchRange = $d.Range($target.Range.Start, $pEnd)\n    foreach ($tok in $tokens) {\n        $find = $searchRange.Find\n        $find.ClearFormatting()\n        $find.Text = $tok\n        $find.Forward = $true\n        $find.Wrap = 0  # wdFindStop: never leave the search range\n        $found = $find.Execute()\n        if (-not $found) {\n            throw \"Token not found: $tok (in paragraph: $paragraphMarker)\"\n        }\n        $searchRange.Font.Bold = 1\n        $searchRange.Font.Color = $HighlightColor\n        # Advance the search window to just after the match so the next\n        # token (or a repeat of the same token) is found in left-to-right\n        # order within this paragraph only.\n        $searchRange.Collapse(0)  # wdCollapseEnd\n        $searchRange.End = $pEnd\n    }\n}\n\n# 1. Data Science & Political Analytics bullets (Siege Analytics)\nHighlight-Tokens 'Discovered systematic race coding errors' @('23%', '64%')\nHighlight-Tokens 'Utilized advanced sampling methods' @([char]0xB1 + '4.2%', [char]0xB1 + '2.1%', '71%', '87%')\nHighlight-Tokens 'Trigonometric algorithm for boundary estimation reduced' @('73.5%', '$4.7M')\nHighlight-Tokens 'Built real-time FEC analysis systems' @('$2')\n\n# 2. Data Products Manager (Helm/Murmuration) bullet\nHighlight-Tokens 'Modernized legacy ETL processes' @('57%')\n\n# 3. KEY ACHIEVEMENTS AND IMPACT bullets\nHighlight-Tokens 'Algorithmic innovation: Pioneered trigonometric boundary estimation' @('73.5%')\nHighlight-Tokens 'savings enabled nonprofit access' @('$4.7M')\nHighlight-Tokens 'Platform impact: Built redistricting system serving' @('12,847')\n"}
